$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOSPITALES")

# ---------------------------------------------------------------------------
# New hospital rows (84-91) - fill in the "constant" columns that every other
# record in the table carries (country code, country name, admin level,
# admin1 type, admin2 type) plus the two brand-new "priorizado" hospitals
# (Mario Catarino Rivas / Leonardo Martinez, both in San Pedro Sula, Cortes)
# and their coordinates, matching the source update.
# ---------------------------------------------------------------------------

$rows = 84..91
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "HND"
    $ws.Range("C$r").Value = "Honduras"
    $ws.Range("D$r").Value = 3
    $ws.Range("H$r").Value = "Departamento"
    $ws.Range("L$r").Value = "Municipio"
}

# Rows 88-89: new hospitals located in San Pedro Sula, Cortes
$ws.Range("G88").Value = "Cortés"
$ws.Range("K88").Value = "San Pedro Sula"
$ws.Range("U88").Value = "Mario Catarino Rivas "
$ws.Range("V88").Value = 15.5245187
$ws.Range("W88").Value = -88.0436936

$ws.Range("G89").Value = "Cortés"
$ws.Range("K89").Value = "San Pedro Sula"
$ws.Range("U89").Value = "Leonardo Martínez"
$ws.Range("V89").Value = 15.5005359
$ws.Range("W89").Value = -88.0310652

# Row 90: coordinates for an already-listed priorizado hospital
$ws.Range("V90").Value = 13.3077843
$ws.Range("W90").Value = -87.2075229

# ---------------------------------------------------------------------------
# Grow the HOSPITALES_HN table / AutoFilter range down to row 93 (it now
# reserves space through the newly appended records) and keep the hidden
# _FilterDatabase defined name lined up with the table.
# ---------------------------------------------------------------------------

$tbl = $ws.ListObjects.Item("HOSPITALES_HN")
$tbl.Resize($ws.Range("A1:W93"))

$filterDbName = $wb.Names.Item("HOSPITALES!_FilterDatabase")
$filterDbName.RefersTo = "=HOSPITALES!`$A`$1:`$W`$93"

# ---------------------------------------------------------------------------
# Match the author's final selection/view state.
# ---------------------------------------------------------------------------

$ws.Range("V91").Select()
